$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "PEP_ID-2007800"
$ws.Range("A8").Value = "PEP_ID-2007802"
